$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 422, shifting existing rows 422:478 down to 423:479
$ws.Rows("422:422").Insert()

# Populate the newly inserted row with the new weekly price entry
$ws.Range("A422").Value = 10
$ws.Range("B422").Value = "Vega Modelo de Temuco"
$ws.Range("C422").Value = "La Araucanía"
$ws.Range("D422").Value2 = 45142
$ws.Range("E422").Value = 9
$ws.Range("F422").Value = "Fruta"
$ws.Range("G422").Value = 100102
$ws.Range("H422").Value = "Cítricos"
$ws.Range("I422").Value = 100102006
$ws.Range("J422").Value = "Pomelo"
$ws.Range("K422").Value = "Start Ruby"
$ws.Range("L422").Value = "Primera"
$ws.Range("M422").Value = 45
$ws.Range("N422").Value = 15000
$ws.Range("O422").Value = 15000
$ws.Range("P422").Value = 15000
$ws.Range("Q422").Value = "$/bandeja 15 kilos granel"
$ws.Range("R422").Value = "Región de O'Higgins"
$ws.Range("S422").Value = 1000
$ws.Range("T422").Value = 15
